# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" sheet (copied from the existing "2022-Q2" sheet so that
#    formatting/column layout matches) positioned right after "总计" and before "2022-Q2".
# 2) Fill it in with the 2022-Q3 fund holdings data.
# 3) Update the "总计" (summary) sheet: insert a new row right after the header with the
#    2022-Q3 totals, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet by duplicating "2022-Q2" (sheet index 2)
# so it starts out with identical styles/column widths, then trim it down to the
# 3 data rows we need and overwrite the values.
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The copied sheet has 9 data rows (rows 2-10); we only need 3 (rows 2-4)
$q3Sheet.Range("A5:H10").Delete()

# Force text storage (so values like "001643", leading zeros, and trailing zeros
# such as "0.9800" are preserved exactly instead of being coerced to numbers).
$q3Sheet.Range("B2:G4").NumberFormat = "@"

$q3Sheet.Range("B2").Value = "001643"
$q3Sheet.Range("C2").Value = "汇丰晋信智造先锋股票A"
$q3Sheet.Range("D2").Value = "15.68"
$q3Sheet.Range("E2").Value = "94.47"
$q3Sheet.Range("F2").Value = "6.25"
$q3Sheet.Range("G2").Value = "0.9800"
$q3Sheet.Range("H2").Value = 3

$q3Sheet.Range("B3").Value = "001644"
$q3Sheet.Range("C3").Value = "汇丰晋信智造先锋股票C"
$q3Sheet.Range("D3").Value = "8.27"
$q3Sheet.Range("E3").Value = "94.47"
$q3Sheet.Range("F3").Value = "6.25"
$q3Sheet.Range("G3").Value = "0.5169"
$q3Sheet.Range("H3").Value = 3

$q3Sheet.Range("B4").Value = "217021"
$q3Sheet.Range("C4").Value = "招商优势企业混合"
$q3Sheet.Range("D4").Value = "2.91"
$q3Sheet.Range("E4").Value = "79.74"
$q3Sheet.Range("F4").Value = "5.08"
$q3Sheet.Range("G4").Value = "0.1478"
$q3Sheet.Range("H4").Value = 8

# ---------------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet - insert a row for 2022-Q3 right
# after the header row, shifting the other quarters down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Rows.Item(2).Insert()

# Copy formatting from row 3 (the row that just got pushed down, which still has
# the original look) onto the freshly inserted, blank row 2.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.64

# Column A is just a sequential 0-based row index, not shifted data - Insert()
# carried the old numbers down along with the rest of the row, so restore the
# correct sequence 0,1,2,3,4 for rows 2-6.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
